$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells keep their original Text format so numeric-looking
# strings (e.g. "311.00", "1.010") are not auto-converted to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.866.01"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "1.828.74"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "311.00"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "0.4577"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "0.3670"
$ws.Range("D9").Value = "0.07168"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "0.8734"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "0.07805"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "19.44"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").Value = "1.871.13"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "5.319"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").Value = "87.17"
$ws.Range("E16").Value = "  -4.71%  "
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "0.000008705"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "26.890.69"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "14.43"
$ws.Range("D22").Value = "4.978"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").Value = "10.45"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "1.996"
$ws.Range("E24").Value = "  +4.14%  "
$ws.Range("D25").Value = "151.62"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "18.13"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "113.63"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "4.904"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").Value = "0.08776"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "3.104"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "0.7387"
$ws.Range("E32").Value = "  -4.27%  "
$ws.Range("D33").Value = "4.473"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "1.128"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").Value = "2.485"
$ws.Range("E35").Value = "  -6.68%  "
$ws.Range("D36").Value = "1.082"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "0.01936"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.913"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05121"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "6.917"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "0.4947"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("D42").Value = "0.1588"
$ws.Range("E42").Value = "  -2.68%  "
$ws.Range("D43").Value = "8.225"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Value = "1.008"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").Value = "0.4644"
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "10.13"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").Value = "102.97"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "1.591"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").Value = "64.58"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("E51").Value = "  -0.91%  "
